$wb = $excel.ActiveWorkbook

# --- Insert the new sheet "Đơn phụ phẫu 1" right after "Đơn sale chính"
# (i.e. right before "Lương"), matching the position in the target file.
$firstSheet = $wb.Worksheets.Item(1)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $firstSheet)
$newSheet.Name = "Đơn phụ phẫu 1"

# --- Header row
$newSheet.Range("A1").Value = "Tiền tố"
$newSheet.Range("B1").Value = "Mã dịch vụ"
$newSheet.Range("C1").Value = "Ngày thực hiện"
$newSheet.Range("D1").Value = "Cơ sở"
$newSheet.Range("E1").Value = "Khách hàng"
$newSheet.Range("F1").Value = "Nguồn khách"
$newSheet.Range("G1").Value = "Tên dịch vụ"
$newSheet.Range("H1").Value = "Phụ phẫu 1"
$newSheet.Range("I1").Value = "Công phụ phẫu 1"

# --- Data row
$newSheet.Range("A2").Value = "HD-LUXURY"
$newSheet.Range("B2").Value = 619
$newSheet.Range("C2").Value = "'08-02-2024"
$newSheet.Range("D2").Value = "SÓC TRĂNG"
$newSheet.Range("E2").Value = "mai hồng nương"
$newSheet.Range("F2").Value = "Cá nhân"
$newSheet.Range("G2").Value = "Thu cánh mũi"
$newSheet.Range("H2").Value = "Kha Như Huỳnh "
$newSheet.Range("I2").Value = 100000

# --- Totals row
$newSheet.Range("A3").Value = "Tổng"
$newSheet.Range("B3").Value = 1
$newSheet.Range("I3").Value = 100000

# --- Update the "Lương" sheet values (commission now flows from the new sheet)
$luong = $wb.Worksheets.Item("Lương")
$luong.Range("B29").Value = 100000
$luong.Range("B34").Value = 100000
$luong.Range("A35").Value = "Tổng lương tại HỆ THỐNG"
$luong.Range("B35").Value = 100000

# Restore original sheet selection/active tab (keep view state close to original)
$firstSheet.Activate()
